$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I is currently empty (data lives in column J). Deleting column I
# shifts column J (and everything to its right) one column to the left,
# so the google_translate data ends up in column I.
$ws.Columns("I").Delete()

# Select the whole column I, matching the selection state after the
# column delete operation.
$ws.Range("I1:I1048576").Select()
